# Regenerate database with new stat order: PTS/REB/AST/STL/BLK
# Swap the values stored in columns I (blocks) and J (steals) for every
# data row in the sheet (rows 2 through 230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 230

for ($r = 2; $r -le $lastRow; $r++) {
    $iCell = $ws.Cells.Item($r, 9)
    $jCell = $ws.Cells.Item($r, 10)

    $iVal = $iCell.Value2
    $jVal = $jCell.Value2

    $iCell.Value = $jVal
    $jCell.Value = $iVal
}
